$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 24, shifting existing rows 24:36 down to 25:37
$ws.Rows("24:24").Insert()

# Populate the new row 24 with the latest weekly data point
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 44755
$ws.Range("D24").NumberFormat = $ws.Range("D25").NumberFormat
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 100112013
$ws.Range("G24").Value = "Alcachofa"
$ws.Range("H24").Value = "Argentina(o)"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16500
$ws.Range("N24").Value = "$/caja 50 unidades"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 330
$ws.Range("Q24").Value = 50
$ws.Range("R24").Value = "Hortaliza"
